# Mise à jour des données au 8 janvier 2018
# - Profil_etudiant : suppression de la ligne "2D" / "UPS D2E"
# - Situation_sociale : suppression de la ligne "DD" / "Demi Droit"

$wb = $excel.ActiveWorkbook

$wsProfil = $wb.Worksheets.Item("Profil_etudiant")
$wsProfil.Rows("3:3").Delete()

$wsSituation = $wb.Worksheets.Item("Situation_sociale")
$wsSituation.Rows("6:6").Delete()
